$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): force text format so numeric-looking strings
# (e.g. "7.200", "0.00001311") are preserved exactly as text, not
# auto-converted/truncated to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.549.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.693.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3932"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3994"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.200"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001311"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.596"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.689.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07061"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.875"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.546.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.061"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.338"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.570"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.874.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.083"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.297"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08540"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.934"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02743"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09023"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.473"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7657"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7150"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.517"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.202"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.333"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07988"
$ws.Range("D51").Style = "Normal"

# Volume(1h) column (E)
$ws.Range("E2").Value = "  +3.17%  "
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("E9").Value = "  +4.81%  "
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("E13").Value = "  +6.57%  "
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("E16").Value = "  +4.16%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("E19").Value = "  +3.92%  "
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("E21").Value = "  +3.37%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("E25").Value = "  +7.78%  "
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("E30").Value = "  +3.34%  "
$ws.Range("E31").Value = "  +12.14%  "
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("E33").Value = "  -3.37%  "
$ws.Range("E34").Value = "  +9.87%  "
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("E36").Value = "  +7.67%  "
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("E40").Value = "  +8.93%  "
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("E46").Value = "  +3.45%  "
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("E49").Value = "  +9.18%  "
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("E51").Value = "  +2.72%  "
